$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Candidates with events")

# Insert a new column at C for the candidate's middle name (shifts
# last_name and everything after it one column to the right).
$ws.Columns("C:C").Insert()

$ws.Range("C1").Value = "candidate_sheet.middle_name"
$ws.Range("C2").Value = "Anne"
$ws.Range("C3").Value = "Richard"
$ws.Range("C4").Value = "baz"
